# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values on the active sheet for rows 2-41
# with freshly recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 3
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 5
    23 = 0
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 2
    30 = 0
    31 = 1
    32 = 2
    33 = 1
    34 = 2
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 3
    41 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
